# Scheduled runner update: refresh profitability figures (currentAveragePrice /
# NQ / HQ price & profit columns, H:N) for a handful of leve rows across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with newly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 181
$ws.Range("I11").Value = 181
$ws.Range("K11").Value = 181
$ws.Range("M11").Value = -41

$ws.Range("H86").Value = 1017.4
$ws.Range("I86").Value = 784.8
$ws.Range("K86").Value = 784.8
$ws.Range("M86").Value = 338.2

$ws.Range("H89").Value = 1017.4
$ws.Range("I89").Value = 784.8
$ws.Range("K89").Value = 3924
$ws.Range("M89").Value = 1692

$ws.Range("H100").Value = 1070.4286
$ws.Range("I100").Value = 1062.25
$ws.Range("J100").Value = 1081.3334
$ws.Range("K100").Value = 1062.25
$ws.Range("L100").Value = 1081.3334
$ws.Range("M100").Value = -521.25
$ws.Range("N100").Value = -2163.3334

$ws.Range("H135").Value = 2032
$ws.Range("I135").Value = 2032
$ws.Range("K135").Value = 18288
$ws.Range("M135").Value = -15753

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H97").Value = 1917.65
$ws.Range("I97").Value = 1789.5714
$ws.Range("J97").Value = 2216.5
$ws.Range("K97").Value = 1789.5714
$ws.Range("L97").Value = 2216.5
$ws.Range("M97").Value = -1293.5714
$ws.Range("N97").Value = -3208.5

$ws.Range("H102").Value = 2875
$ws.Range("I102").Value = 2750
$ws.Range("K102").Value = 2750
$ws.Range("M102").Value = -1128

$ws.Range("H104").Value = 20357.143
$ws.Range("J104").Value = 20357.143
$ws.Range("L104").Value = 20357.143
$ws.Range("N104").Value = -27345.143

$ws.Range("H132").Value = 5482.4
$ws.Range("I132").Value = 4678
$ws.Range("J132").Value = 8700
$ws.Range("K132").Value = 14034
$ws.Range("L132").Value = 26100
$ws.Range("M132").Value = -11504
$ws.Range("N132").Value = -31160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 14948
$ws.Range("I33").Value = 20021
$ws.Range("J33").Value = 9875
$ws.Range("K33").Value = 20021
$ws.Range("L33").Value = 9875
$ws.Range("M33").Value = -19685
$ws.Range("N33").Value = -10547

$ws.Range("H86").Value = 2262.8
$ws.Range("I86").Value = 2291
$ws.Range("J86").Value = 2150
$ws.Range("K86").Value = 2291
$ws.Range("L86").Value = 2150
$ws.Range("M86").Value = -1168
$ws.Range("N86").Value = -4396

$ws.Range("H89").Value = 2262.8
$ws.Range("I89").Value = 2291
$ws.Range("J89").Value = 2150
$ws.Range("K89").Value = 11455
$ws.Range("L89").Value = 10750
$ws.Range("M89").Value = -5839
$ws.Range("N89").Value = -21982

$ws.Range("H94").Value = 2413.7856
$ws.Range("I94").Value = 1974.25
$ws.Range("J94").Value = 2999.8333
$ws.Range("K94").Value = 1974.25
$ws.Range("L94").Value = 2999.8333
$ws.Range("M94").Value = -1523.25
$ws.Range("N94").Value = -3901.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2075
$ws.Range("I16").Value = 1150
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1150
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -863
$ws.Range("N16").Value = -3574

$ws.Range("H36").Value = 7511.75
$ws.Range("I36").Value = 7511.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 7511.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -7123.75
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 7511.75
$ws.Range("I40").Value = 7511.75
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7511.75
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7351.75
$ws.Range("N40").ClearContents()

$ws.Range("H113").Value = 2075
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 450
$ws.Range("I50").Value = 450
$ws.Range("K50").Value = 1350
$ws.Range("M50").Value = -869

$ws.Range("H53").Value = 450
$ws.Range("I53").Value = 450
$ws.Range("K53").Value = 1350
$ws.Range("M53").Value = -869

$ws.Range("H60").Value = 244.5
$ws.Range("I60").Value = 174.5
$ws.Range("K60").Value = 523.5
$ws.Range("M60").Value = -272.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 25500
$ws.Range("J47").Value = 25500
$ws.Range("L47").Value = 25500
$ws.Range("N47").Value = -26636

$ws.Range("H97").Value = 377444.62
$ws.Range("I97").Value = 2333.8
$ws.Range("J97").Value = 1002629.3
$ws.Range("K97").Value = 2333.8
$ws.Range("L97").Value = 1002629.3
$ws.Range("M97").Value = -1837.8
$ws.Range("N97").Value = -1003621.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2080.5
$ws.Range("I16").Value = 1949.1428
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1949.1428
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1779.1428
$ws.Range("N16").Value = -3340

$ws.Range("H32").Value = 5907.2
$ws.Range("I32").Value = 1134
$ws.Range("K32").Value = 1134
$ws.Range("M32").Value = -817

$ws.Range("H46").Value = 2146.25
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2834
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2834
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -3210

$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -798

$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 10000
$ws.Range("N71").Value = -17488

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2995
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4238

$ws.Range("H65").Value = 2995
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21190

$ws.Range("H96").Value = 1919
$ws.Range("I96").Value = 1925
$ws.Range("J96").Value = 1913
$ws.Range("K96").Value = 1925
$ws.Range("L96").Value = 1913
$ws.Range("M96").Value = -552
$ws.Range("N96").Value = -4659

$ws.Range("H126").Value = 1388
$ws.Range("I126").Value = 1388
$ws.Range("K126").Value = 4164
$ws.Range("M126").Value = -1694
